$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (prefixed with a literal quote
# character below to force Excel to treat numeric-looking strings as text,
# matching the original inlineStr/shared-string cell type).
$updates = @{
    "D2" = '27.899.28'
    "E2" = '  +0.83%  '
    "D3" = '1.858.80'
    "E3" = '  +0.32%  '
    "D4" = '0.9994'
    "E4" = '  -0.12%  '
    "D5" = '316.60'
    "E5" = '  +1.24%  '
    "D6" = '1.000'
    "E6" = '  -0.02%  '
    "D7" = '0.4324'
    "E7" = '  +1.25%  '
    "D8" = '0.3706'
    "E8" = '  +1.61%  '
    "D9" = '45.36'
    "E9" = '  +1.65%  '
    "D10" = '0.07383'
    "E10" = '  +1.22%  '
    "D11" = '0.8814'
    "E11" = '  +0.75%  '
    "D12" = '21.22'
    "E12" = '  +2.91%  '
    "D13" = '1.918.27'
    "E13" = '  +5.22%  '
    "D14" = '5.467'
    "E14" = '  +2.94%  '
    "D15" = '6.625'
    "E15" = '  +1.69%  '
    "D16" = '0.06972'
    "E16" = '  +1.01%  '
    "D17" = '1.003'
    "E17" = '  +0.11%  '
    "D18" = '81.72'
    "E18" = '  +2.41%  '
    "D19" = '0.000009092'
    "E19" = '  +0.94%  '
    "D20" = '1.0000'
    "E20" = '  +0.04%  '
    "D21" = '15.64'
    "E21" = '  +2.10%  '
    "D22" = '27.945.73'
    "E22" = '  +0.93%  '
    "D23" = '5.095'
    "E23" = '  +2.52%  '
    "E24" = '  +6.18%  '
    "D25" = '2.084.29'
    "E25" = '  +1.15%  '
    "D26" = '1.965'
    "E26" = '  +0.14%  '
    "D27" = '155.08'
    "E27" = '  +1.25%  '
    "D28" = '18.93'
    "E28" = '  +0.47%  '
    "D29" = '5.338'
    "E29" = '  +1.37%  '
    "D30" = '115.94'
    "E30" = '  -5.09%  '
    "D31" = '1.863'
    "E31" = '  -0.32%  '
    "D32" = '0.08946'
    "D33" = '0.7906'
    "E33" = '  +3.16%  '
    "D34" = '4.631'
    "E34" = '  +2.35%  '
    "B35" = 'ARBITRUM'
    "C35" = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    "D35" = '1.183'
    "E35" = '  +7.07%  '
    "B36" = 'HuobiToken'
    "C36" = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    "D36" = '2.972'
    "E36" = '  +0.20%  '
    "D37" = '0.9994'
    "E37" = '  -0.01%  '
    "B38" = 'Hedera'
    "C38" = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    "D38" = '0.05461'
    "E38" = '  +1.50%  '
    "B39" = 'TrustWalletToken'
    "C39" = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    "D39" = '1.109'
    "E39" = '  +1.56%  '
    "D40" = '0.01965'
    "E40" = '  +1.66%  '
    "D41" = '2.843'
    "E41" = '  +0.69%  '
    "D42" = '0.5202'
    "E42" = '  +2.87%  '
    "D43" = '0.1687'
    "E43" = '  +2.37%  '
    "D44" = '6.793'
    "E44" = '  -0.03%  '
    "D45" = '8.726'
    "E45" = '  +4.52%  '
    "D46" = '10.65'
    "E46" = '  +3.29%  '
    "D47" = '0.4808'
    "E47" = '  +3.20%  '
    "D48" = '107.05'
    "E48" = '  +2.20%  '
    "D49" = '0.06577'
    "E49" = '  +0.63%  '
    "B50" = 'PaxDollar'
    "C50" = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    "D50" = '1.0000'
    "E50" = '  +0.07%  '
    "B51" = 'NEARProtocol'
    "C51" = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    "D51" = '1.670'
    "E51" = '  +2.80%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Leading apostrophe forces Excel to store the value as text, preserving
    # exact formatting (leading/trailing zeros, thousand-dot separators, etc.)
    $cell.Value = "'" + $updates[$ref]
    $cell.Style = "Normal"
}
